# Auto-generated edit: populate market/profit columns (H:N) for
# specific leve rows on the BSM and LTW sheets, matching the
# scheduled-runner market data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BSM")

# Rows with no realised market activity: H:L all zero.
$zeroRows_BSM = @(117, 118, 119, 120, 122, 123, 124, 125, 126, 127, 128, 129, 130, 131, 132, 133, 135, 138, 139, 140, 141)
foreach ($r in $zeroRows_BSM) {
    foreach ($c in 8..12) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Row 134: realised leve turn-in with recorded market prices.
$ws.Range("H134").Value = 4999
$ws.Range("I134").Value = 4999
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14997
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12462

# Row 137: realised leve turn-in with recorded market prices.
$ws.Range("H137").Value = 100000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 100000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200


$ws = $wb.Worksheets.Item("LTW")

# Rows with no realised market activity: H:L all zero.
$zeroRows_LTW = @(124, 125, 128, 129, 130, 131, 133, 134, 135, 137, 138, 139, 140, 141)
foreach ($r in $zeroRows_LTW) {
    foreach ($c in 8..12) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Row 127: realised leve turn-in with recorded market prices.
$ws.Range("H127").Value = 117000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 117000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 117000
$ws.Range("N127").Value = -126920

# Row 132: realised leve turn-in with recorded market prices.
$ws.Range("H132").Value = 3517.111
$ws.Range("I132").Value = 1951.2
$ws.Range("J132").Value = 5474.5
$ws.Range("K132").Value = 5853.6
$ws.Range("L132").Value = 16423.5
$ws.Range("M132").Value = -3323.6
$ws.Range("N132").Value = -21483.5

# Row 136: realised leve turn-in with recorded market prices.
$ws.Range("H136").Value = 2668.8333
$ws.Range("I136").Value = 1003.3333
$ws.Range("J136").Value = 4334.3335
$ws.Range("K136").Value = 3009.9999
$ws.Range("L136").Value = 13003.0005
$ws.Range("M136").Value = -459.9998999999998
$ws.Range("N136").Value = -18103.0005

